$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 329, shifting existing rows 329-409 down to 330-410.
$ws.Rows.Item(329).Insert()

# Populate the new row 329 with the new data record.
$ws.Cells.Item(329, 1).Value = 4
$ws.Cells.Item(329, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(329, 3).Value = "Los Lagos"
$ws.Cells.Item(329, 4).Value = 44932
$ws.Cells.Item(329, 5).Value = 10
$ws.Cells.Item(329, 6).Value = 100112045
$ws.Cells.Item(329, 7).Value = "Zapallo"
$ws.Cells.Item(329, 8).Value = "Paine"
$ws.Cells.Item(329, 9).Value = "1a nueva(o)"
$ws.Cells.Item(329, 10).Value = 1200
$ws.Cells.Item(329, 11).Value = 600
$ws.Cells.Item(329, 12).Value = 650
$ws.Cells.Item(329, 13).Value = 625
$ws.Cells.Item(329, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(329, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(329, 16).Value = 625
$ws.Cells.Item(329, 17).Value = 1
$ws.Cells.Item(329, 18).Value = "Hortaliza"

# Ensure date column style/format for the new row matches the rest (column D date format).
$ws.Range("D329").NumberFormat = "YYYY-MM-DD HH:MM:SS"
